$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("učitelja Mat. ", $true, $false, $false, $false, $false, $true, 1, $false, "{{ radno_mj }} ", 2)
